# Generate Report for Handoff
# A new handoff event occurred for b.md (zh-cn and de-de) on 2016-08-12.
# Update the Overview sheet and the per-locale status sheets to reflect
# the freshly generated handoff xliff files and the resulting "not latest"
# handback warning.

$wb = $excel.ActiveWorkbook

# ---- Overview sheet : row for b.md (row 3) ----
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E3").Value = "Ready for handoff"
$overview.Range("F3").Value = "Ready for handoff"
$overview.Range("G3").Value = "2016-08-12 14:42:24"

# ---- zh-cn sheet : row for b.md (row 3) ----
# Leading apostrophe forces these to stay plain text (shared strings)
# instead of Excel auto-coercing "False"/"True"-looking text into a
# Boolean cell type.
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = "Ready for handoff"
$zhcn.Range("F3").Value = "'False"
$zhcn.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$zhcn.Range("H3").Value = "2016-08-12 14:42:17"
$zhcn.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/oltest/blob/d13721fe50643e6c7efeedfd4bca28911179c3b0/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/oltest/blob/d1eda731b9ac61a24d3561c932367ecf7a5870d2/e2e/b.md."

# ---- de-de sheet : row for b.md (row 3) ----
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = "Ready for handoff"
$dede.Range("F3").Value = "'False"
$dede.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$dede.Range("H3").Value = "2016-08-12 14:42:24"
$dede.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/oltest/blob/d13721fe50643e6c7efeedfd4bca28911179c3b0/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/oltest/blob/d1eda731b9ac61a24d3561c932367ecf7a5870d2/e2e/b.md."

# The Error Detail column (P) on both locale sheets now holds a long
# message, so widen column P the same way Excel's autofit would.
# (ColumnWidth is specified in characters; the saved OOXML width ends up
# 0.8333... wider than the ColumnWidth value due to the standard column
# padding, so back that constant out to land on an OOXML width of 40.)
$colWidthPad = 0.8333333333333334
$zhcn.Columns.Item(16).ColumnWidth = 40 - $colWidthPad
$dede.Columns.Item(16).ColumnWidth = 40 - $colWidthPad
